$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AV3").Value2 = -0.06834955435453159
$ws.Range("AZ3").Value2 = 0
$ws.Range("BA3").Value2 = -0.01048984217249318
$ws.Range("BE3").Value2 = -0.06346346529129036
$ws.Range("BI3").Value2 = 0
$ws.Range("BJ3").Value2 = -0.01350133615972497
$ws.Range("BN3").Value2 = -0.06914825245171299
$ws.Range("BR3").Value2 = 0
$ws.Range("BS3").Value2 = -0.01149420626780502
$ws.Range("CB3").Value2 = -0.0008498064264525762
$ws.Range("CX3").Value2 = -0.02589212009742492
$ws.Range("CY3").Value2 = -0.007585994751205218
$ws.Range("CZ3").Value2 = 0
$ws.Range("L5").Value2 = 281336.94137175
$ws.Range("AP5").Value2 = 0.103173098434489
$ws.Range("AQ5").Value2 = 0.1620134241373689
$ws.Range("AR5").Value2 = 0.1326312370746145
$ws.Range("AS5").Value2 = 0.04202660262247519
$ws.Range("AT5").Value2 = 0.02598690690414279
$ws.Range("AU5").Value2 = 0.004546482354745358
$ws.Range("AV5").Value2 = 0.0561659616600887
$ws.Range("AW5").Value2 = 0.03246848306892767
$ws.Range("AX5").Value2 = 0.07147750213957864
$ws.Range("AY5").Value2 = 0.02435248903068241
$ws.Range("AZ5").Value2 = 0.0428438280422095
$ws.Range("BA5").Value2 = 0.02127519230028399
$ws.Range("BB5").Value2 = 0.03994782635701525
$ws.Range("BC5").Value2 = 0.02649496179872661
$ws.Range("BD5").Value2 = 0.006428064423120565
$ws.Range("BE5").Value2 = 0.0536987082513803
$ws.Range("BF5").Value2 = 0.03323001182290512
$ws.Range("BG5").Value2 = 0.06866222409854451
$ws.Range("BH5").Value2 = 0.02470944814757399
$ws.Range("BI5").Value2 = 0.04277051884462714
$ws.Range("BJ5").Value2 = 0.02207955549375783
$ws.Range("BK5").Value2 = 0.04119123823787647
$ws.Range("BL5").Value2 = 0.02525954850203326
$ws.Range("BM5").Value2 = 0.004334972861881025
$ws.Range("BN5").Value2 = 0.05574637722610173
$ws.Range("BO5").Value2 = 0.03281560254581525
$ws.Range("BP5").Value2 = 0.07107233374979242
$ws.Range("BQ5").Value2 = 0.0246067353789282
$ws.Range("BR5").Value2 = 0.04297030232282382
$ws.Range("BS5").Value2 = 0.02184627511891423
$ws.Range("BV5").Value2 = 0.04387625478466212
$ws.Range("BX5").Value2 = 0.01098508208894688
$ws.Range("BY5").Value2 = 1.551183090186825
$ws.Range("BZ5").Value2 = 0.2251971447103978
$ws.Range("CA5").Value2 = 1.537213581130476
$ws.Range("CB5").Value2 = 0.1609030712433077
$ws.Range("CX5").Value2 = 0.01649981738486565
$ws.Range("CY5").Value2 = 0.01232115011941043
$ws.Range("CZ5").Value2 = 0.0610565643784701
$ws.Range("DA5").Value2 = 0.3470901431481323
$ws.Range("DB5").Value2 = 0.3970634181197694
$ws.Range("DC5").Value2 = 1.283851118541529
$ws.Range("L6").Value2 = 535627.0019095155
$ws.Range("AP6").Value2 = 0.4969604033047024
$ws.Range("AQ6").Value2 = 0.5155465927595724
$ws.Range("AR6").Value2 = 0.5130559177684323
$ws.Range("AS6").Value2 = 0.09904420824242778
$ws.Range("AT6").Value2 = 0.07778123016219282
$ws.Range("AU6").Value2 = 0.9803599790578759
$ws.Range("AV6").Value2 = 0.05449563262391398
$ws.Range("AW6").Value2 = 0.1396392754666966
$ws.Range("AX6").Value2 = 0.2038560643924456
$ws.Range("AY6").Value2 = 0.05803675588273972
$ws.Range("AZ6").Value2 = 0.08881311245039997
$ws.Range("BA6").Value2 = 0.02165919305182338
$ws.Range("BB6").Value2 = 0.09899013204910595
$ws.Range("BC6").Value2 = 0.07781667858776187
$ws.Range("BD6").Value2 = 0.9790350228377188
$ws.Range("BE6").Value2 = 0.05232357591105089
$ws.Range("BF6").Value2 = 0.1347167893108235
$ws.Range("BG6").Value2 = 0.2081997699272716
$ws.Range("BH6").Value2 = 0.05751976539683133
$ws.Range("BI6").Value2 = 0.08759615952415677
$ws.Range("BJ6").Value2 = 0.02117567457511895
$ws.Range("BK6").Value2 = 0.1015886278086519
$ws.Range("BL6").Value2 = 0.07888757077532256
$ws.Range("BM6").Value2 = 0.9807784471780124
$ws.Range("BN6").Value2 = 0.05536661078171171
$ws.Range("BO6").Value2 = 0.1396596698386499
$ws.Range("BP6").Value2 = 0.2050253433611077
$ws.Range("BQ6").Value2 = 0.05883633615191687
$ws.Range("BR6").Value2 = 0.09157084500895245
$ws.Range("BS6").Value2 = 0.02177042348978891
$ws.Range("BV6").Value2 = 0.08590109189150916
$ws.Range("BX6").Value2 = 0.009609004169864481
$ws.Range("BY6").Value2 = 1.306417433279428
$ws.Range("BZ6").Value2 = 0.1262993104407133
$ws.Range("CA6").Value2 = 1.289392838358581
$ws.Range("CB6").Value2 = 0.08722257781536563
$ws.Range("CX6").Value2 = 0.009051893742401791
$ws.Range("CY6").Value2 = 0.006976538998277805
$ws.Range("CZ6").Value2 = 0.0361121418681964
$ws.Range("DA6").Value2 = 0.2312730809530294
$ws.Range("DB6").Value2 = 0.291522785768328
$ws.Range("DC6").Value2 = 1.902287474159097
$ws.Range("AP7").Value2 = 0.4955623721707529
$ws.Range("AR7").Value2 = 0.5130360317720301
$ws.Range("AS7").Value2 = 0.08385812800130335
$ws.Range("AT7").Value2 = 0.07221450703691004
$ws.Range("AU7").Value2 = 0.9808384293699252
$ws.Range("AV7").Value2 = 0.0324771139014379
$ws.Range("AW7").Value2 = 0.1379
$ws.Range("AX7").Value2 = 0.1739708279247994
$ws.Range("AY7").Value2 = 0.05234270711431947
$ws.Range("AZ7").Value2 = 0.07229327079430101
$ws.Range("BA7").Value2 = 0.01256745275660578
$ws.Range("BB7").Value2 = 0.08742236473856652
$ws.Range("BC7").Value2 = 0.0729217898938159
$ws.Range("BD7").Value2 = 0.9807362704447131
$ws.Range("BE7").Value2 = 0.03202242077805297
$ws.Range("BF7").Value2 = 0.1327474664384867
$ws.Range("BG7").Value2 = 0.184857666396061
$ws.Range("BH7").Value2 = 0.05188404032611363
$ws.Range("BI7").Value2 = 0.07143272014100209
$ws.Range("BJ7").Value2 = 0.01194563815547413
$ws.Range("BK7").Value2 = 0.08829458792488301
$ws.Range("BL7").Value2 = 0.07486981573233652
$ws.Range("BM7").Value2 = 0.9811220535802468
$ws.Range("BN7").Value2 = 0.03365596965556104
$ws.Range("BO7").Value2 = 0.1379
$ws.Range("BP7").Value2 = 0.1764406694175532
$ws.Range("BQ7").Value2 = 0.05327858031580367
$ws.Range("BR7").Value2 = 0.07585494157496497
$ws.Range("BS7").Value2 = 0.01240064718688126
$ws.Range("BV7").Value2 = 0.06894036018385516
$ws.Range("DC7").Value2 = 1.800110952098913
$ws.Range("L8").Value2 = 129500298202.5096
$ws.Range("AP8").Value2 = 0.02128181719081568
$ws.Range("AQ8").Value2 = 0.03828590951330329
$ws.Range("AR8").Value2 = 0.03163243250163777
$ws.Range("AS8").Value2 = 0.005364212887409809
$ws.Range("AT8").Value2 = 0.002291191358480511
$ws.Range("AU8").Value2 = 0.0003908150210877011
$ws.Range("AV8").Value2 = 0.009143146773307144
$ws.Range("AW8").Value2 = 0.003754259483012278
$ws.Range("AX8").Value2 = 0.01209640701351525
$ws.Range("AY8").Value2 = 0.002503449475820893
$ws.Range("AZ8").Value2 = 0.006086096109210615
$ws.Range("BA8").Value2 = 0.002707541416401497
$ws.Range("BB8").Value2 = 0.005215403409720813
$ws.Range("BC8").Value2 = 0.00246951855259184
$ws.Range("BD8").Value2 = 0.001971500763427782
$ws.Range("BE8").Value2 = 0.008784988176392033
$ws.Range("BF8").Value2 = 0.003858425684357852
$ws.Range("BG8").Value2 = 0.01179196432732534
$ws.Range("BH8").Value2 = 0.002615620401613208
$ws.Range("BI8").Value2 = 0.006386861710499212
$ws.Range("BJ8").Value2 = 0.003032183417102996
$ws.Range("BK8").Value2 = 0.005375313772722254
$ws.Range("BL8").Value2 = 0.002270937091724466
$ws.Range("BM8").Value2 = 0.0002847780332430891
$ws.Range("BN8").Value2 = 0.009112873498313591
$ws.Range("BO8").Value2 = 0.003794356635281106
$ws.Range("BP8").Value2 = 0.01207902652541699
$ws.Range("BQ8").Value2 = 0.002580082053191583
$ws.Range("BR8").Value2 = 0.006268185437733821
$ws.Range("BS8").Value2 = 0.00283617294121849
$ws.Range("BV8").Value2 = 0.006174640654042887
$ws.Range("BX8").Value2 = 0.004970142988546125
$ws.Range("BY8").Value2 = 4.811606754842686
$ws.Range("BZ8").Value2 = 0.1563874698473492
$ws.Range("CA8").Value2 = 4.724581178153104
$ws.Range("CB8").Value2 = 0.1036160498354396
$ws.Range("CX8").Value2 = 0.007954677340277362
$ws.Range("CY8").Value2 = 0.01198404209740771
$ws.Range("CZ8").Value2 = 0.03597030699899115
$ws.Range("DA8").Value2 = 0.4844132396024792
$ws.Range("DB8").Value2 = 0.361275899296724
$ws.Range("DC8").Value2 = 3.108681222546991
$ws.Range("L9").Value2 = 359861.4986387258
$ws.Range("AP9").Value2 = 0.1458828886155456
$ws.Range("AQ9").Value2 = 0.1956678550843323
$ws.Range("AR9").Value2 = 0.1778550884895841
$ws.Range("AS9").Value2 = 0.07324078704799539
$ws.Range("AT9").Value2 = 0.0478663906982813
$ws.Range("AU9").Value2 = 0.01976904198709945
$ws.Range("AV9").Value2 = 0.09561980324863226
$ws.Range("AW9").Value2 = 0.06127201223244001
$ws.Range("AX9").Value2 = 0.1099836670306789
$ws.Range("AY9").Value2 = 0.0500344828675274
$ws.Range("AZ9").Value2 = 0.07801343544038176
$ws.Range("BA9").Value2 = 0.05203404093861534
$ws.Range("BB9").Value2 = 0.0722177499630168
$ws.Range("BC9").Value2 = 0.04969425069957127
$ws.Range("BD9").Value2 = 0.04440158514544026
$ws.Range("BE9").Value2 = 0.09372826775520837
$ws.Range("BF9").Value2 = 0.06211622722250485
$ws.Range("BG9").Value2 = 0.1085908114313791
$ws.Range("BH9").Value2 = 0.05114313640766675
$ws.Range("BI9").Value2 = 0.07991784350505969
$ws.Range("BJ9").Value2 = 0.05506526506885258
$ws.Range("BK9").Value2 = 0.07331653137405134
$ws.Range("BL9").Value2 = 0.04765435018678218
$ws.Range("BM9").Value2 = 0.0168753676476422
$ws.Range("BN9").Value2 = 0.09546137176006633
$ws.Range("BO9").Value2 = 0.06159834929022941
$ws.Range("BP9").Value2 = 0.1099046246771126
$ws.Range("BQ9").Value2 = 0.05079450810069513
$ws.Range("BR9").Value2 = 0.07917187276889326
$ws.Range("BS9").Value2 = 0.05325573153397192
$ws.Range("BV9").Value2 = 0.07857888173067168
$ws.Range("BX9").Value2 = 0.07049924104943348
$ws.Range("BY9").Value2 = 2.193537497934031
$ws.Range("BZ9").Value2 = 0.3954585564219709
$ws.Range("CA9").Value2 = 2.17361017161613
$ws.Range("CB9").Value2 = 0.3218944700292932
$ws.Range("CX9").Value2 = 0.08918899786564126
$ws.Range("CY9").Value2 = 0.1094716497427882
$ws.Range("CZ9").Value2 = 0.1896583955404852
$ws.Range("DA9").Value2 = 0.6959980169529789
$ws.Range("DB9").Value2 = 0.6010623089969326
$ws.Range("DC9").Value2 = 1.763145264164865
